$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: several Price (column D) values look like plain decimal numbers
# (e.g. "1.00", "160.00", "8.00") and must be forced to remain literal text
# -- matching the inline-string cell type used throughout the sheet -- rather
# than being auto-converted to a number (which would also drop trailing zeros).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "74.405.44"
$ws.Range("E2").Value = "  +8.29%  "
$ws.Range("D3").Value = "2.596.95"
$ws.Range("E3").Value = "  +6.99%  "
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.01%  "
Set-TextValue $ws.Range("D5") "186.87"
$ws.Range("E5").Value = "  +15.70%  "
Set-TextValue $ws.Range("D6") "587.19"
$ws.Range("E6").Value = "  +5.03%  "
Set-TextValue $ws.Range("D7") "1.00"
$ws.Range("E7").Value = "  -0.09%  "
Set-TextValue $ws.Range("D8") "0.538"
$ws.Range("E8").Value = "  +5.03%  "
Set-TextValue $ws.Range("D9") "0.206"
$ws.Range("E9").Value = "  +23.00%  "
$ws.Range("D10").Value = "2.596.14"
$ws.Range("E10").Value = "  +7.09%  "
$ws.Range("E11").Value = "  +0.13%  "
Set-TextValue $ws.Range("D12") "0.363"
$ws.Range("E12").Value = "  +10.00%  "
Set-TextValue $ws.Range("D13") "4.80"
$ws.Range("E13").Value = "  +4.33%  "
$ws.Range("E14").Value = "  +9.62%  "
$ws.Range("D15").Value = "74.514.81"
$ws.Range("E15").Value = "  +8.58%  "
$ws.Range("D16").Value = "3.072.90"
$ws.Range("E16").Value = "  +6.83%  "
Set-TextValue $ws.Range("D17") "26.35"
$ws.Range("E17").Value = "  +14.00%  "
$ws.Range("D18").Value = "2.598.89"
$ws.Range("E18").Value = "  +6.91%  "
Set-TextValue $ws.Range("D19") "9.15"
$ws.Range("E19").Value = "  +32.40%  "
Set-TextValue $ws.Range("D20") "11.83"
$ws.Range("E20").Value = "  +13.17%  "
Set-TextValue $ws.Range("D21") "375.66"
$ws.Range("E21").Value = "  +11.68%  "
Set-TextValue $ws.Range("D22") "2.29"
$ws.Range("E22").Value = "  +18.79%  "
Set-TextValue $ws.Range("D23") "4.10"
$ws.Range("E23").Value = "  +7.68%  "
$ws.Range("E24").Value = "  +0.07%  "
Set-TextValue $ws.Range("D25") "70.27"
$ws.Range("E25").Value = "  +5.11%  "
Set-TextValue $ws.Range("D26") "4.18"
$ws.Range("E26").Value = "  +13.57%  "
Set-TextValue $ws.Range("D27") "9.35"
$ws.Range("E27").Value = "  +14.62%  "
$ws.Range("D28").Value = "2.730.69"
$ws.Range("E28").Value = "  +6.85%  "
Set-TextValue $ws.Range("D29") "1.00"
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("E30").Value = "  +16.77%  "
Set-TextValue $ws.Range("D31") "1.41"
$ws.Range("E31").Value = "  +22.51%  "
Set-TextValue $ws.Range("D32") "8.00"
$ws.Range("E32").Value = "  +12.52%  "
Set-TextValue $ws.Range("D33") "509.14"
$ws.Range("E33").Value = "  +19.40%  "
Set-TextValue $ws.Range("D34") "1.76"
$ws.Range("E34").Value = "  +9.22%  "
Set-TextValue $ws.Range("D35") "0.999"
$ws.Range("E35").Value = "  -0.04%  "
Set-TextValue $ws.Range("D36") "0.122"
$ws.Range("E36").Value = "  +14.56%  "
Set-TextValue $ws.Range("D37") "160.00"
$ws.Range("E37").Value = "  +0.19%  "
Set-TextValue $ws.Range("D38") "19.23"
$ws.Range("E38").Value = "  +7.46%  "
Set-TextValue $ws.Range("D39") "19.36"
$ws.Range("E39").Value = "  +1.80%  "
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("E41").Value = "  +13.82%  "
$ws.Range("E42").Value = "  +13.21%  "
Set-TextValue $ws.Range("D43") "0.328"
$ws.Range("E43").Value = "  +10.41%  "
Set-TextValue $ws.Range("D44") "2.41"
$ws.Range("E44").Value = "  +18.83%  "
Set-TextValue $ws.Range("D45") "156.83"
$ws.Range("E45").Value = "  +19.76%  "
Set-TextValue $ws.Range("D48") "38.80"
$ws.Range("E48").Value = "  +3.52%  "
Set-TextValue $ws.Range("D49") "3.64"
$ws.Range("E49").Value = "  +9.06%  "
Set-TextValue $ws.Range("D50") "0.525"
$ws.Range("E50").Value = "  +9.49%  "
Set-TextValue $ws.Range("D51") "20.40"
$ws.Range("E51").Value = "  +21.47%  "

# Rows 46 & 47: the two coins swap places (Coin name + Link), and each
# picks up new Price / Volume(1h) figures.
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D46") "0.0870"
$ws.Range("E46").Value = "  +21.77%  "

$ws.Range("B47").Value = "ImmutableX"
$ws.Range("C47").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D47") "1.18"
$ws.Range("E47").Value = "  +9.60%  "
